# Generate Report for Handoff
# Updates the "Latest Handoff"/"Latest HO Xliff Generate Date" timestamps for the
# 0452522f-50ef-4efc-92f1-477dcd1e679a file (row 5 in every sheet), reflecting a
# freshly generated handoff xliff for both the zh-cn and de-de locales, and the
# overview's rolled-up "Latest HO Xliff Generate Date" (max across locales).

$wb = $excel.ActiveWorkbook

# zh-cn sheet: "Latest Handoff Datetime" (column H) for row 5 (0452522f...)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H5").Value = "2016-09-01 04:46:13"

# de-de sheet: "Latest Handoff Datetime" (column H) for row 5 (0452522f...)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").Value = "2016-09-01 04:46:18"

# Overview sheet: "Latest HO Xliff Generate Date" (column G) for row 5 (0452522f...)
# reflects the newest of the per-locale handoff datetimes above.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G5").Value = "2016-09-01 04:46:18"
